$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H5").Value = 417.83334
$ws.Range("I5").Value = 101.4
$ws.Range("K5").Value = 101.4
$ws.Range("M5").Value = 13.59999999999999
$ws.Range("H11").Value = 827.8571
$ws.Range("I11").Value = 827.8571
$ws.Range("K11").Value = 827.8571
$ws.Range("M11").Value = -687.8571
$ws.Range("H12").Value = 7462.7856
$ws.Range("I12").Value = 8011.4614
$ws.Range("J12").Value = 330
$ws.Range("K12").Value = 8011.4614
$ws.Range("L12").Value = 330
$ws.Range("M12").Value = -7841.4614
$ws.Range("N12").Value = -670
$ws.Range("H29").Value = 2048.2727
$ws.Range("J29").Value = 3556.8333
$ws.Range("L29").Value = 10670.4999
$ws.Range("N29").Value = -11232.4999
$ws.Range("H33").Value = 534.375
$ws.Range("J33").Value = 565.8
$ws.Range("L33").Value = 565.8
$ws.Range("N33").Value = -1023.8
$ws.Range("H38").Value = 1024.1
$ws.Range("I38").Value = 1024.1
$ws.Range("K38").Value = 3072.3
$ws.Range("M38").Value = -2700.3
$ws.Range("H39").Value = 1573.3125
$ws.Range("I39").Value = 1166
$ws.Range("J39").Value = 1980.625
$ws.Range("K39").Value = 3498
$ws.Range("L39").Value = 5941.875
$ws.Range("M39").Value = -3202
$ws.Range("N39").Value = -6533.875
$ws.Range("H41").Value = 1165.5238
$ws.Range("I41").Value = 1385.3334
$ws.Range("K41").Value = 1385.3334
$ws.Range("M41").Value = -945.3334
$ws.Range("H64").Value = 100002250
$ws.Range("J64").Value = 0
$ws.Range("L64").Value = 0
$ws.Range("N64").ClearContents()
$ws.Range("H67").Value = 100002250
$ws.Range("J67").Value = 0
$ws.Range("L67").Value = 0
$ws.Range("N67").ClearContents()
$ws.Range("H70").Value = 12534.538
$ws.Range("I70").Value = 2012.25
$ws.Range("K70").Value = 6036.75
$ws.Range("M70").Value = -5766.75
$ws.Range("H73").Value = 12534.538
$ws.Range("I73").Value = 2012.25
$ws.Range("K73").Value = 6036.75
$ws.Range("M73").Value = -5100.75
$ws.Range("H82").Value = 5619.6665
$ws.Range("I82").Value = 4521.143
$ws.Range("K82").Value = 13563.429
$ws.Range("M82").Value = -13157.429
$ws.Range("H85").Value = 5619.6665
$ws.Range("I85").Value = 4521.143
$ws.Range("K85").Value = 13563.429
$ws.Range("M85").Value = -12159.429
$ws.Range("H98").Value = 1945.3704
$ws.Range("I98").Value = 1794.7916
$ws.Range("K98").Value = 1794.7916
$ws.Range("M98").Value = -296.7916
$ws.Range("H100").Value = 827.5
$ws.Range("J100").Value = 953.1429
$ws.Range("L100").Value = 953.1429
$ws.Range("N100").Value = -2035.1429
$ws.Range("H106").Value = 5000
$ws.Range("J106").Value = 5000
$ws.Range("L106").Value = 5000
$ws.Range("N106").Value = -6262
$ws.Range("H122").Value = 1945.3704
$ws.Range("I122").Value = 1794.7916
$ws.Range("K122").Value = 5384.3748
$ws.Range("M122").Value = -2934.3748
$ws.Range("H131").Value = 2324.8125
$ws.Range("I131").Value = 1302.3572
$ws.Range("K131").Value = 3907.0716
$ws.Range("M131").Value = 1132.9284
$ws.Range("H132").Value = 3441.6667
$ws.Range("I132").Value = 3441.6667
$ws.Range("K132").Value = 10325.0001
$ws.Range("M132").Value = -7795.000100000001
$ws.Range("H135").Value = 22727772
$ws.Range("I135").Value = 23810026
$ws.Range("K135").Value = 214290234
$ws.Range("M135").Value = -214287699
$ws.Range("H137").Value = 2202.9778
$ws.Range("I137").Value = 1956.2122
$ws.Range("K137").Value = 5868.6366
$ws.Range("M137").Value = -3318.6366
$ws.Range("H138").Value = 2626.6724
$ws.Range("J138").Value = 3758.0344
$ws.Range("L138").Value = 11274.1032
$ws.Range("N138").Value = -21554.1032
$ws.Range("H141").Value = 957.7037
$ws.Range("I141").Value = 957.7037
$ws.Range("K141").Value = 2873.1111
$ws.Range("M141").Value = 2306.8889

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3042.2693
$ws.Range("I32").Value = 3076.5476
$ws.Range("J32").Value = 2898.3
$ws.Range("K32").Value = 3076.5476
$ws.Range("L32").Value = 2898.3
$ws.Range("M32").Value = -2789.5476
$ws.Range("N32").Value = -3472.3
$ws.Range("H45").Value = 3344.375
$ws.Range("I45").Value = 3556
$ws.Range("J45").Value = 2991.6667
$ws.Range("K45").Value = 3556
$ws.Range("L45").Value = 2991.6667
$ws.Range("M45").Value = -3179
$ws.Range("N45").Value = -3745.6667
$ws.Range("H61").Value = 62501704
$ws.Range("I61").Value = 76924680
$ws.Range("K61").Value = 76924680
$ws.Range("M61").Value = -76924468
$ws.Range("H63").Value = 3039.7334
$ws.Range("I63").Value = 3039.7334
$ws.Range("K63").Value = 3039.7334
$ws.Range("M63").Value = -2353.7334
$ws.Range("H64").Value = 0
$ws.Range("I64").Value = 0
$ws.Range("K64").Value = 0
$ws.Range("M64").ClearContents()
$ws.Range("H66").Value = 3039.7334
$ws.Range("I66").Value = 3039.7334
$ws.Range("K66").Value = 15198.667
$ws.Range("M66").Value = -11766.667
$ws.Range("H67").Value = 0
$ws.Range("I67").Value = 0
$ws.Range("K67").Value = 0
$ws.Range("M67").ClearContents()
$ws.Range("H74").Value = 47627470
$ws.Range("I74").Value = 55563440
$ws.Range("J74").Value = 11666.667
$ws.Range("K74").Value = 55563440
$ws.Range("L74").Value = 11666.667
$ws.Range("M74").Value = -55562566
$ws.Range("N74").Value = -13414.667
$ws.Range("H77").Value = 47627470
$ws.Range("I77").Value = 55563440
$ws.Range("J77").Value = 11666.667
$ws.Range("K77").Value = 277817200
$ws.Range("L77").Value = 58333.335
$ws.Range("M77").Value = -277812832
$ws.Range("N77").Value = -67069.33499999999
$ws.Range("H102").Value = 50000250
$ws.Range("I102").Value = 50000250
$ws.Range("J102").Value = 0
$ws.Range("K102").Value = 50000250
$ws.Range("L102").Value = 0
$ws.Range("M102").Value = -49998628
$ws.Range("N102").ClearContents()
$ws.Range("H110").Value = 49654.855
$ws.Range("I110").Value = 64022.375
$ws.Range("K110").Value = 64022.375
$ws.Range("M110").Value = -61977.375
$ws.Range("H132").Value = 3128146.2
$ws.Range("I132").Value = 3706791.2
$ws.Range("J132").Value = 3463
$ws.Range("K132").Value = 11120373.6
$ws.Range("L132").Value = 10389
$ws.Range("M132").Value = -11117843.6
$ws.Range("N132").Value = -15449
$ws.Range("H136").Value = 62501704
$ws.Range("I136").Value = 76924680
$ws.Range("K136").Value = 230774040
$ws.Range("M136").Value = -230771490

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H62").Value = 94285.71
$ws.Range("I62").Value = 60000
$ws.Range("J62").Value = 95555.555
$ws.Range("K62").Value = 60000
$ws.Range("L62").Value = 95555.555
$ws.Range("M62").Value = -59314
$ws.Range("N62").Value = -96927.555
$ws.Range("H65").Value = 94285.71
$ws.Range("I65").Value = 60000
$ws.Range("J65").Value = 95555.555
$ws.Range("K65").Value = 180000
$ws.Range("L65").Value = 286666.665
$ws.Range("M65").Value = -176568
$ws.Range("N65").Value = -293530.665
$ws.Range("H86").Value = 3644.4443
$ws.Range("I86").Value = 3725
$ws.Range("J86").Value = 3000
$ws.Range("K86").Value = 3725
$ws.Range("L86").Value = 3000
$ws.Range("M86").Value = -2602
$ws.Range("N86").Value = -5246
$ws.Range("H89").Value = 3644.4443
$ws.Range("I89").Value = 3725
$ws.Range("J89").Value = 3000
$ws.Range("K89").Value = 18625
$ws.Range("L89").Value = 15000
$ws.Range("M89").Value = -13009
$ws.Range("N89").Value = -26232
$ws.Range("H92").Value = 37999
$ws.Range("J92").Value = 37999
$ws.Range("L92").Value = 37999
$ws.Range("N92").Value = -42991
$ws.Range("H94").Value = 13844.546
$ws.Range("I94").Value = 13839.1
$ws.Range("K94").Value = 13839.1
$ws.Range("M94").Value = -13388.1
$ws.Range("H99").Value = 2092.7896
$ws.Range("I99").Value = 1994.375
$ws.Range("J99").Value = 2164.3635
$ws.Range("K99").Value = 1994.375
$ws.Range("L99").Value = 2164.3635
$ws.Range("M99").Value = -496.375
$ws.Range("N99").Value = -5160.363499999999
$ws.Range("H103").Value = 24999.75
$ws.Range("J103").Value = 24999.75
$ws.Range("L103").Value = 24999.75
$ws.Range("N103").Value = -27343.75
$ws.Range("H134").Value = 57223148
$ws.Range("I134").Value = 57223148
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 171669444
$ws.Range("L134").Value = 0
$ws.Range("M134").Value = -171666909
$ws.Range("N134").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 2482.1914
$ws.Range("I22").Value = 2583.432
$ws.Range("J22").Value = 997.3333
$ws.Range("K22").Value = 2583.432
$ws.Range("L22").Value = 997.3333
$ws.Range("M22").Value = -2233.432
$ws.Range("N22").Value = -1697.3333
$ws.Range("H31").Value = 9439.161
$ws.Range("I31").Value = 8325.667
$ws.Range("K31").Value = 8325.667
$ws.Range("M31").Value = -8030.666999999999
$ws.Range("H34").Value = 9439.161
$ws.Range("I34").Value = 8325.667
$ws.Range("K34").Value = 8325.667
$ws.Range("M34").Value = -8123.666999999999
$ws.Range("H51").Value = 18666.5
$ws.Range("J51").Value = 18400
$ws.Range("L51").Value = 18400
$ws.Range("N51").Value = -19872
$ws.Range("H59").Value = 71989
$ws.Range("I59").Value = 13580.2
$ws.Range("J59").Value = 145000
$ws.Range("K59").Value = 13580.2
$ws.Range("L59").Value = 145000
$ws.Range("M59").Value = -12435.2
$ws.Range("N59").Value = -147290
$ws.Range("H61").Value = 18666.5
$ws.Range("J61").Value = 18400
$ws.Range("L61").Value = 18400
$ws.Range("N61").Value = -19096
$ws.Range("H62").Value = 5899.6665
$ws.Range("I62").Value = 5799
$ws.Range("K62").Value = 5799
$ws.Range("M62").Value = -5175
$ws.Range("H65").Value = 5899.6665
$ws.Range("I65").Value = 5799
$ws.Range("K65").Value = 28995
$ws.Range("M65").Value = -25875
$ws.Range("H86").Value = 5086.364
$ws.Range("I86").Value = 4805.778
$ws.Range("J86").Value = 6349
$ws.Range("K86").Value = 4805.778
$ws.Range("L86").Value = 6349
$ws.Range("M86").Value = -3682.778
$ws.Range("N86").Value = -8595
$ws.Range("H89").Value = 5086.364
$ws.Range("I89").Value = 4805.778
$ws.Range("J89").Value = 6349
$ws.Range("K89").Value = 24028.89
$ws.Range("L89").Value = 31745
$ws.Range("M89").Value = -18412.89
$ws.Range("N89").Value = -42977
$ws.Range("H94").Value = 2192
$ws.Range("I94").Value = 1538.4
$ws.Range("K94").Value = 1538.4
$ws.Range("M94").Value = -1087.4
$ws.Range("H95").Value = 18208
$ws.Range("J95").Value = 18208
$ws.Range("L95").Value = 18208
$ws.Range("N95").Value = -23700
$ws.Range("H96").Value = 15000
$ws.Range("J96").Value = 15000
$ws.Range("L96").Value = 15000
$ws.Range("N96").Value = -20492
$ws.Range("H99").Value = 3190.25
$ws.Range("I99").Value = 3190.25
$ws.Range("K99").Value = 3190.25
$ws.Range("M99").Value = -1692.25
$ws.Range("H122").Value = 2923.926
$ws.Range("I122").Value = 2923.926
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 8771.778
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -6321.778
$ws.Range("N122").ClearContents()
$ws.Range("H126").Value = 3190.25
$ws.Range("I126").Value = 3190.25
$ws.Range("K126").Value = 9570.75
$ws.Range("M126").Value = -7100.75
$ws.Range("H132").Value = 34485204
$ws.Range("I132").Value = 38464070
$ws.Range("J132").Value = 1684.6666
$ws.Range("K132").Value = 115392210
$ws.Range("L132").Value = 5053.9998
$ws.Range("M132").Value = -115389680
$ws.Range("N132").Value = -10113.9998
$ws.Range("H134").Value = 5103504.5
$ws.Range("I134").Value = 5815356.5
$ws.Range("J134").Value = 1899.5
$ws.Range("K134").Value = 17446069.5
$ws.Range("L134").Value = 5698.5
$ws.Range("M134").Value = -17443534.5
$ws.Range("N134").Value = -10768.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 10676.667
$ws.Range("I3").Value = 10676.667
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 32030.001
$ws.Range("L3").Value = 0
$ws.Range("M3").Value = -31918.001
$ws.Range("N3").ClearContents()
$ws.Range("H5").Value = 251449.5
$ws.Range("I5").Value = 500399
$ws.Range("K5").Value = 1501197
$ws.Range("M5").Value = -1501085
$ws.Range("H12").Value = 193.73334
$ws.Range("J12").Value = 157.625
$ws.Range("L12").Value = 472.875
$ws.Range("N12").Value = -818.875
$ws.Range("H32").Value = 1534.8
$ws.Range("J32").Value = 1534.8
$ws.Range("L32").Value = 4604.4
$ws.Range("N32").Value = -5170.4
$ws.Range("H33").Value = 2105.4285
$ws.Range("J33").Value = 2164.6667
$ws.Range("L33").Value = 12988.0002
$ws.Range("N33").Value = -13554.0002
$ws.Range("H37").Value = 149974.6
$ws.Range("J37").Value = 149974.6
$ws.Range("L37").Value = 449923.8
$ws.Range("N37").Value = -450147.8
$ws.Range("H60").Value = 5970
$ws.Range("I60").Value = 595
$ws.Range("K60").Value = 1785
$ws.Range("M60").Value = -1534
$ws.Range("H68").Value = 4445.75
$ws.Range("I68").Value = 4461.3335
$ws.Range("K68").Value = 13384.0005
$ws.Range("M68").Value = -12573.0005
$ws.Range("H71").Value = 4445.75
$ws.Range("I71").Value = 4461.3335
$ws.Range("K71").Value = 40152.0015
$ws.Range("M71").Value = -36096.0015
$ws.Range("H113").Value = 143621
$ws.Range("I113").Value = 500374.5
$ws.Range("J113").Value = 919.6
$ws.Range("K113").Value = 1501123.5
$ws.Range("L113").Value = 2758.8
$ws.Range("M113").Value = -1498953.5
$ws.Range("N113").Value = -7098.8
$ws.Range("H122").Value = 380.30435
$ws.Range("J122").Value = 585.6667
$ws.Range("L122").Value = 5271.0003
$ws.Range("N122").Value = -10171.0003
$ws.Range("H123").Value = 5923.7144
$ws.Range("I123").Value = 0
$ws.Range("J123").Value = 5923.7144
$ws.Range("K123").Value = 0
$ws.Range("L123").Value = 17771.1432
$ws.Range("M123").ClearContents()
$ws.Range("N123").Value = -22671.1432
$ws.Range("H131").Value = 1388.375
$ws.Range("J131").Value = 4900
$ws.Range("L131").Value = 14700
$ws.Range("N131").Value = -24780
$ws.Range("H135").Value = 251449.5
$ws.Range("I135").Value = 500399
$ws.Range("K135").Value = 4503591
$ws.Range("M135").Value = -4501056
$ws.Range("H137").Value = 12501474
$ws.Range("H139").Value = 1385
$ws.Range("I139").Value = 1385
$ws.Range("K139").Value = 4155
$ws.Range("M139").Value = 985

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H64").Value = 60000
$ws.Range("J64").Value = 60000
$ws.Range("L64").Value = 60000
$ws.Range("N64").Value = -60496
$ws.Range("H67").Value = 60000
$ws.Range("J67").Value = 60000
$ws.Range("L67").Value = 60000
$ws.Range("N67").Value = -61716
$ws.Range("H80").Value = 3234.5
$ws.Range("I80").Value = 2368.4443
$ws.Range("K80").Value = 2368.4443
$ws.Range("M80").Value = -1370.4443
$ws.Range("H83").Value = 3234.5
$ws.Range("I83").Value = 2368.4443
$ws.Range("K83").Value = 11842.2215
$ws.Range("M83").Value = -6850.2215
$ws.Range("H95").Value = 0
$ws.Range("J95").Value = 0
$ws.Range("L95").Value = 0
$ws.Range("N95").ClearContents()
$ws.Range("H107").Value = 1737.1
$ws.Range("I107").Value = 1820.3334
$ws.Range("K107").Value = 1820.3334
$ws.Range("M107").Value = 99.66660000000002
$ws.Range("H122").Value = 4516.514
$ws.Range("I122").Value = 2679.625
$ws.Range("J122").Value = 8524.272
$ws.Range("K122").Value = 8038.875
$ws.Range("L122").Value = 25572.816
$ws.Range("M122").Value = -5588.875
$ws.Range("N122").Value = -30472.816
$ws.Range("H132").Value = 3677523.2
$ws.Range("I132").Value = 4167745.8
$ws.Range("J132").Value = 855.25
$ws.Range("K132").Value = 12503237.4
$ws.Range("L132").Value = 2565.75
$ws.Range("M132").Value = -12500707.4
$ws.Range("N132").Value = -7625.75
$ws.Range("H135").Value = 86664.336
$ws.Range("I135").Value = 0
$ws.Range("J135").Value = 86664.336
$ws.Range("K135").Value = 0
$ws.Range("L135").Value = 86664.336
$ws.Range("M135").ClearContents()
$ws.Range("N135").Value = -96804.336

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 22749
$ws.Range("J2").Value = 22749
$ws.Range("L2").Value = 22749
$ws.Range("N2").Value = -22973
$ws.Range("H16").Value = 2209.7368
$ws.Range("I16").Value = 1201.9375
$ws.Range("J16").Value = 2942.682
$ws.Range("K16").Value = 1201.9375
$ws.Range("L16").Value = 2942.682
$ws.Range("M16").Value = -1031.9375
$ws.Range("N16").Value = -3282.682
$ws.Range("H22").Value = 3147.0715
$ws.Range("I22").Value = 3137.375
$ws.Range("J22").Value = 3160
$ws.Range("K22").Value = 3137.375
$ws.Range("L22").Value = 3160
$ws.Range("M22").Value = -2842.375
$ws.Range("N22").Value = -3750
$ws.Range("H27").Value = 3147.0715
$ws.Range("I27").Value = 3137.375
$ws.Range("J27").Value = 3160
$ws.Range("K27").Value = 3137.375
$ws.Range("L27").Value = 3160
$ws.Range("M27").Value = -3030.375
$ws.Range("N27").Value = -3374
$ws.Range("H61").Value = 3085
$ws.Range("I61").Value = 3166.6667
$ws.Range("K61").Value = 3166.6667
$ws.Range("M61").Value = -2964.6667
$ws.Range("H68").Value = 8671666
$ws.Range("J68").Value = 999999
$ws.Range("L68").Value = 999999
$ws.Range("N68").Value = -1001497
$ws.Range("H71").Value = 8671666
$ws.Range("J71").Value = 999999
$ws.Range("L71").Value = 4999995
$ws.Range("N71").Value = -5007483
$ws.Range("H82").Value = 1524.7142
$ws.Range("I82").Value = 1505.75
$ws.Range("J82").Value = 1550
$ws.Range("K82").Value = 1505.75
$ws.Range("L82").Value = 1550
$ws.Range("M82").Value = -1144.75
$ws.Range("N82").Value = -2272
$ws.Range("H85").Value = 1524.7142
$ws.Range("I85").Value = 1505.75
$ws.Range("J85").Value = 1550
$ws.Range("K85").Value = 1505.75
$ws.Range("L85").Value = 1550
$ws.Range("M85").Value = -257.75
$ws.Range("N85").Value = -4046
$ws.Range("H100").Value = 15356802
$ws.Range("I100").Value = 19962344
$ws.Range("K100").Value = 19962344
$ws.Range("M100").Value = -19961803
$ws.Range("H102").Value = 69999
$ws.Range("J102").Value = 69999
$ws.Range("L102").Value = 69999
$ws.Range("N102").Value = -76489
$ws.Range("H113").Value = 3085
$ws.Range("I113").Value = 3166.6667
$ws.Range("K113").Value = 3166.6667
$ws.Range("M113").Value = -996.6667000000002
$ws.Range("H132").Value = 14493608
$ws.Range("I132").Value = 16140053
$ws.Range("J132").Value = 4899.4
$ws.Range("K132").Value = 48420159
$ws.Range("L132").Value = 14698.2
$ws.Range("M132").Value = -48417629
$ws.Range("N132").Value = -19758.2
$ws.Range("H136").Value = 1733.8182
$ws.Range("I136").Value = 905.5
$ws.Range("J136").Value = 2424.0833
$ws.Range("K136").Value = 2716.5
$ws.Range("L136").Value = 7272.249899999999
$ws.Range("M136").Value = -166.5
$ws.Range("N136").Value = -12372.2499

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H5").Value = 29999
$ws.Range("J5").Value = 29999
$ws.Range("L5").Value = 29999
$ws.Range("N5").Value = -30223
$ws.Range("H64").Value = 27100
$ws.Range("I64").Value = 27100
$ws.Range("J64").Value = 0
$ws.Range("K64").Value = 27100
$ws.Range("L64").Value = 0
$ws.Range("M64").Value = -26852
$ws.Range("N64").ClearContents()
$ws.Range("H67").Value = 27100
$ws.Range("I67").Value = 27100
$ws.Range("J67").Value = 0
$ws.Range("K67").Value = 27100
$ws.Range("L67").Value = 0
$ws.Range("M67").Value = -26242
$ws.Range("N67").ClearContents()
$ws.Range("H81").Value = 251758.75
$ws.Range("I81").Value = 251758.75
$ws.Range("J81").Value = 0
$ws.Range("K81").Value = 503517.5
$ws.Range("L81").Value = 0
$ws.Range("M81").Value = -502456.5
$ws.Range("N81").ClearContents()
$ws.Range("H84").Value = 251758.75
$ws.Range("I84").Value = 251758.75
$ws.Range("J84").Value = 0
$ws.Range("K84").Value = 2517587.5
$ws.Range("L84").Value = 0
$ws.Range("M84").Value = -2512283.5
$ws.Range("N84").ClearContents()
$ws.Range("H122").Value = 1225.0625
$ws.Range("I122").Value = 1225.0625
$ws.Range("K122").Value = 3675.1875
$ws.Range("M122").Value = -1225.1875
$ws.Range("H126").Value = 1736.125
$ws.Range("I126").Value = 1347.5
$ws.Range("K126").Value = 4042.5
$ws.Range("M126").Value = -1572.5
$ws.Range("H132").Value = 15156570
$ws.Range("I132").Value = 19232730
$ws.Range("J132").Value = 16550.715
$ws.Range("K132").Value = 57698190
$ws.Range("L132").Value = 49652.145
$ws.Range("M132").Value = -57695660
$ws.Range("N132").Value = -54712.145
$ws.Range("H136").Value = 11629684
$ws.Range("I136").Value = 12501876
$ws.Range("J136").Value = 467.66666
$ws.Range("K136").Value = 37505628
$ws.Range("L136").Value = 1402.99998
$ws.Range("M136").Value = -37503078
$ws.Range("N136").Value = -6502.999980000001
